# Replace the division-fact answers in the table cells.
# Source text strings repeat ("63÷2=31, 1" and "10÷6=1, 4" each occur twice),
# so cells are addressed positionally (row, column) via the Tables API
# instead of a blind global Find/Replace.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$newValues = @(
    @("20÷8=2, 4", "88÷9=9, 7", "93÷5=18, 3", "72÷6=12, 0", "35÷6=5, 5"),
    @("95÷9=10, 5", "79÷2=39, 1", "13÷8=1, 5", "34÷4=8, 2", "26÷4=6, 2"),
    @("56÷5=11, 1", "62÷9=6, 8", "34÷7=4, 6", "89÷6=14, 5", "80÷7=11, 3"),
    @("83÷6=13, 5", "80÷2=40, 0", "26÷9=2, 8", "80÷6=13, 2", "39÷9=4, 3"),
    @("65÷6=10, 5", "19÷4=4, 3", "12÷5=2, 2", "97÷4=24, 1", "30÷6=5, 0")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowIndex = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($rowIndex, $c)
        $rng = $cell.Range
        # Trim the trailing cell-mark / paragraph-mark characters so only
        # the visible text is replaced, preserving the run's formatting.
        $rng.End = $rng.End - 1
        $rng.Text = $newValues[$r][$c - 1]
    }
}

$d.Save()
